$d = $word.ActiveDocument

# --- Table 1: Bill header info (Date / ID No / Name / TP No) ---
$t1 = $d.Tables.Item(1)
$t1.Cell(1,1).Range.Text = "දිනය"
$t1.Cell(1,2).Range.Text = "2020-12-28 13:04:54"
$t1.Cell(2,1).Range.Text = "ජා. අංකය"
$t1.Cell(2,2).Range.Text = "123"
$t1.Cell(3,1).Range.Text = "නම"
$t1.Cell(3,2).Range.Text = "Test"
$t1.Cell(4,1).Range.Text = "දු. අංකය"
$t1.Cell(4,2).Range.Text = "123"

# --- Table 2: Rented Item details ---
$t2 = $d.Tables.Item(2)

# Header row
$t2.Cell(1,1).Range.Text = "භාණ්ඩ වර්ගය"
$t2.Cell(1,2).Range.Text = "කුලියට ගත් දිනය"
$t2.Cell(1,3).Range.Text = "දින ගණන"
$t2.Cell(1,4).Range.Text = "ප්‍රමාණය"
$t2.Cell(1,5).Range.Text = "දිනකට කුලිය"
$t2.Cell(1,6).Range.Text = "මුදල"

# Row 2 (Poker)
$t2.Cell(2,2).Range.Text = "2020-12-20"
$t2.Cell(2,3).Range.Text = "8.0"
$t2.Cell(2,4).Range.Text = "1"
$t2.Cell(2,6).Range.Text = "16,000.00"

# Row 3 (Concrete Mixer)
$t2.Cell(3,2).Range.Text = "2020-12-20"
$t2.Cell(3,3).Range.Text = "8.0"
$t2.Cell(3,6).Range.Text = "28,000.00"

# Row 4 (Total)
$t2.Cell(4,6).Range.Text = " 44,000.00"

# --- Table 3: Payment details ---
$t3 = $d.Tables.Item(3)
$t3.Cell(1,1).Range.Text = "හිඟ මුදල"
$t3.Cell(2,1).Range.Text = "මෙම බිල්පතෙහි වටිනාකම"
$t3.Cell(2,2).Range.Text = "Rs.  44,000.00"
$t3.Cell(3,1).Range.Text = "මුළු මුදල"
$t3.Cell(3,2).Range.Text = "            Rs.  43,000.00"
$t3.Cell(4,1).Range.Text = "ගෙවීම්"
$t3.Cell(4,2).Range.Text = "-(Rs.  30,000.00)"
$t3.Cell(5,1).Range.Text = "ගෙවිය යුතු වටිනාකම"
$t3.Cell(5,2).Range.Text = "            Rs.  13,000.00"
